# Update NATMI TPM output values (Cthrc1-Fzd3) on Sheet1 rows 2-10
# to reflect the new TPM-based recomputation of ligand/receptor
# expression, specificity, and edge-weight statistics.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.026182
$ws.Range("H2").Value = 0.078546
$ws.Range("I2").Value = 0.02060098669457318
$ws.Range("J2").Value = 0.02060098669457318
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.196431
$ws.Range("N2").Value = 0.589293
$ws.Range("O2").Value = 0.09717285149889213
$ws.Range("P2").Value = 0.09717285149889213
$ws.Range("Q2").Value = 0.005142956442
$ws.Range("R2").Value = 0.046286607978
$ws.Range("S2").Value = 0.002001856620802412
$ws.Range("T2").Value = 0.002001856620802412
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.026182
$ws.Range("H3").Value = 0.078546
$ws.Range("I3").Value = 0.02060098669457318
$ws.Range("J3").Value = 0.02060098669457318
$ws.Range("M3").Value = 0.4307096666666667
$ws.Range("O3").Value = 0.2130686423127578
$ws.Range("P3").Value = 0.2130686423127578
$ws.Range("Q3").Value = 0.01127684049266667
$ws.Range("R3").Value = 0.101491564434
$ws.Range("S3").Value = 0.004389424265315896
$ws.Range("T3").Value = 0.004389424265315895
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.026182
$ws.Range("H4").Value = 0.078546
$ws.Range("I4").Value = 0.02060098669457318
$ws.Range("J4").Value = 0.02060098669457318
$ws.Range("O4").Value = 0.68975850618835
$ws.Range("P4").Value = 0.68975850618835
$ws.Range("Q4").Value = 0.036506060058
$ws.Range("R4").Value = 0.328554540522
$ws.Range("S4").Value = 0.01420970580845487
$ws.Range("T4").Value = 0.01420970580845487
$ws.Range("I5").Value = 0.9231010325934437
$ws.Range("J5").Value = 0.9231010325934434
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.196431
$ws.Range("N5").Value = 0.589293
$ws.Range("O5").Value = 0.09717285149889213
$ws.Range("P5").Value = 0.09717285149889213
$ws.Range("Q5").Value = 0.230448593195
$ws.Range("R5").Value = 2.074037338755
$ws.Range("S5").Value = 0.08970035955867668
$ws.Range("T5").Value = 0.08970035955867667
$ws.Range("I6").Value = 0.9231010325934437
$ws.Range("J6").Value = 0.9231010325934434
$ws.Range("M6").Value = 0.4307096666666667
$ws.Range("O6").Value = 0.2130686423127578
$ws.Range("P6").Value = 0.2130686423127578
$ws.Range("Q6").Value = 0.5052992488905557
$ws.Range("R6").Value = 4.547693240015001
$ws.Range("S6").Value = 0.1966838837321898
$ws.Range("T6").Value = 0.1966838837321898
$ws.Range("I7").Value = 0.9231010325934437
$ws.Range("J7").Value = 0.9231010325934434
$ws.Range("O7").Value = 0.68975850618835
$ws.Range("P7").Value = 0.68975850618835
$ws.Range("S7").Value = 0.6367167893025771
$ws.Range("T7").Value = 0.6367167893025769
$ws.Range("I8").Value = 0.05629798071198328
$ws.Range("J8").Value = 0.05629798071198327
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.196431
$ws.Range("N8").Value = 0.589293
$ws.Range("O8").Value = 0.09717285149889213
$ws.Range("P8").Value = 0.09717285149889213
$ws.Range("Q8").Value = 0.014054572573
$ws.Range("R8").Value = 0.126491153157
$ws.Range("S8").Value = 0.005470635319413044
$ws.Range("T8").Value = 0.005470635319413043
$ws.Range("I9").Value = 0.05629798071198328
$ws.Range("J9").Value = 0.05629798071198327
$ws.Range("M9").Value = 0.4307096666666667
$ws.Range("O9").Value = 0.2130686423127578
$ws.Range("P9").Value = 0.2130686423127578
$ws.Range("S9").Value = 0.0119953343152521
$ws.Range("T9").Value = 0.0119953343152521
$ws.Range("I10").Value = 0.05629798071198328
$ws.Range("J10").Value = 0.05629798071198327
$ws.Range("O10").Value = 0.68975850618835
$ws.Range("P10").Value = 0.68975850618835
$ws.Range("R10").Value = 0.897867537093
$ws.Range("S10").Value = 0.03883201107731812
$ws.Range("T10").Value = 0.03883201107731812
